$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 88 (tezos) and row 89 (gala) data for columns A, B, D, E
# New row 88 = gala data (with updated F value); new row 89 = tezos data (with updated F value)
$ws.Range("A88").Value = "gala"
$ws.Range("B88").Value = "Gala"
$ws.Range("D88").Value = 43744189534.59282
$ws.Range("E88").Value = "https://ethplorer.io/es/address/0x15d4c048f83bd7e37d49ea4c83a07267ec4203da#chart=candlestick"

$ws.Range("A89").Value = "tezos"
$ws.Range("B89").Value = "Tezos"
$ws.Range("D89").Value = 1036827531.678072
$ws.Range("E89").Value = "https://tzkt.io/"

# Update marketCapUsd (column F) values for rows 2-101
$ws.Range("F2").Value = 1645600371530.644
$ws.Range("F3").Value = 221238742768.5913
$ws.Range("F4").Value = 144656039999.9366
$ws.Range("F5").Value = 124997297363.8779
$ws.Range("F6").Value = 87163354463.41571
$ws.Range("F7").Value = 64106576736.51845
$ws.Range("F8").Value = 60326188982.85223
$ws.Range("F9").Value = 25323753239.18353
$ws.Range("F10").Value = 23771977841.44719
$ws.Range("F11").Value = 21875136033.27025
$ws.Range("F12").Value = 17108030022.05006
$ws.Range("F13").Value = 10622166943.44786
$ws.Range("F14").Value = 9118365895.604147
$ws.Range("F15").Value = 8908742037.096777
$ws.Range("F16").Value = 8861972999.71876
$ws.Range("F17").Value = 8236068214.605865
$ws.Range("F18").Value = 8167767101.384493
$ws.Range("F19").Value = 7487771074.991312
$ws.Range("F20").Value = 7319371445.253803
$ws.Range("F21").Value = 7291815763.072726
$ws.Range("F22").Value = 6469254257.014988
$ws.Range("F23").Value = 6211954278.539178
$ws.Range("F24").Value = 6144102402.931064
$ws.Range("F25").Value = 6104196066.729373
$ws.Range("F26").Value = 6023188746.348
$ws.Range("F27").Value = 5577893743.817565
$ws.Range("F28").Value = 5406233793.749472
$ws.Range("F29").Value = 5364422326.337841
$ws.Range("F30").Value = 5245728803.728033
$ws.Range("F31").Value = 4298824719.001237
$ws.Range("F32").Value = 3993800282.934098
$ws.Range("F33").Value = 3769326792.395563
$ws.Range("F34").Value = 3154710197.981006
$ws.Range("F35").Value = 3097008436.659281
$ws.Range("F36").Value = 3029773532.516821
$ws.Range("F37").Value = 2891620175.439909
$ws.Range("F38").Value = 2790157658.150105
$ws.Range("F39").Value = 2699961505.205756
$ws.Range("F40").Value = 2584061240.994352
$ws.Range("F41").Value = 2581685391.143952
$ws.Range("F42").Value = 2549908349.032105
$ws.Range("F43").Value = 2531820394.514572
$ws.Range("F44").Value = 2516755433.877529
$ws.Range("F45").Value = 2109448878.910554
$ws.Range("F46").Value = 2017932366.900606
$ws.Range("F47").Value = 1972372076.705628
$ws.Range("F48").Value = 1968864590.643154
$ws.Range("F49").Value = 1931516772.054979
$ws.Range("F50").Value = 1914050754.452449
$ws.Range("F51").Value = 1861481627.624409
$ws.Range("F52").Value = 1801854633.955202
$ws.Range("F53").Value = 1767555384.779251
$ws.Range("F54").Value = 1667781804.208688
$ws.Range("F55").Value = 1654706339.06789
$ws.Range("F56").Value = 1572698970.098183
$ws.Range("F57").Value = 1541312353.215098
$ws.Range("F58").Value = 1467334025.427559
$ws.Range("F59").Value = 1410827960.660558
$ws.Range("F60").Value = 1369938502.867266
$ws.Range("F61").Value = 1356745023.470731
$ws.Range("F62").Value = 1283915640.470609
$ws.Range("F63").Value = 1246676762.603356
$ws.Range("F64").Value = 1191234770.115511
$ws.Range("F65").Value = 1099987973.711774
$ws.Range("F66").Value = 1097148630.121393
$ws.Range("F67").Value = 1070156987.868036
$ws.Range("F68").Value = 976331197.3122771
$ws.Range("F69").Value = 961915715.5932056
$ws.Range("F70").Value = 951874736.5640842
$ws.Range("F71").Value = 944028010.9010398
$ws.Range("F72").Value = 898582459.0946836
$ws.Range("F73").Value = 885191982.2847631
$ws.Range("F74").Value = 882390168.3221568
$ws.Range("F75").Value = 873887205.0400928
$ws.Range("F76").Value = 868354716.6944952
$ws.Range("F77").Value = 851562953.7912272
$ws.Range("F78").Value = 840759571.291217
$ws.Range("F79").Value = 836039972.3196337
$ws.Range("F80").Value = 830198277.1416909
$ws.Range("F81").Value = 806973804.9338775
$ws.Range("F82").Value = 801516669.7799932
$ws.Range("F83").Value = 795127890.3225529
$ws.Range("F84").Value = 788485386.8367231
$ws.Range("F85").Value = 767629346.5031015
$ws.Range("F86").Value = 696775670.0897274
$ws.Range("F87").Value = 684048671.4149623
$ws.Range("F88").Value = 675934600.2598318
$ws.Range("F89").Value = 673831428.5401824
$ws.Range("F90").Value = 670010329.3811026
$ws.Range("F91").Value = 663315970.4451092
$ws.Range("F92").Value = 656440612.0699849
$ws.Range("F93").Value = 643333898.4821957
$ws.Range("F94").Value = 642893675.9930117
$ws.Range("F95").Value = 618550480.827413
$ws.Range("F96").Value = 612140479.4307131
$ws.Range("F97").Value = 603786260.1301708
$ws.Range("F98").Value = 594222462.900689
$ws.Range("F99").Value = 579262544.1824132
$ws.Range("F100").Value = 572038782.6309732
$ws.Range("F101").Value = 565351361.5889218

Write-Output "done"
